$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "95-77="  # was "12+16="
$t.Cell(1,2).Range.Text = "38-32="  # was "20-17="
$t.Cell(1,3).Range.Text = "33+25="  # was "72-48="
$t.Cell(1,4).Range.Text = "47-26="  # was "53-30="
$t.Cell(1,5).Range.Text = "25+32="  # was "11+62="
$t.Cell(2,1).Range.Text = "63-59="  # was "35+32="
$t.Cell(2,2).Range.Text = "20+26="  # was "40+50="
$t.Cell(2,3).Range.Text = "63-22="  # was "23+61="
$t.Cell(2,4).Range.Text = "19+57="  # was "55-25="
$t.Cell(2,5).Range.Text = "53+29="  # was "91-57="
$t.Cell(3,1).Range.Text = "8+19="  # was "69+26="
$t.Cell(3,2).Range.Text = "67+27="  # was "67+30="
$t.Cell(3,3).Range.Text = "39-19="  # was "19+49="
$t.Cell(3,4).Range.Text = "10+2="  # was "95-94="
$t.Cell(3,5).Range.Text = "89-67="  # was "67+31="
$t.Cell(4,1).Range.Text = "44-41="  # was "19-0="
$t.Cell(4,2).Range.Text = "75-52="  # was "82-32="
$t.Cell(4,3).Range.Text = "66-24="  # was "29+39="
$t.Cell(4,4).Range.Text = "46-33="  # was "95-36="
$t.Cell(4,5).Range.Text = "59-11="  # was "76-36="
$t.Cell(5,1).Range.Text = "47-46="  # was "36+28="
$t.Cell(5,2).Range.Text = "80-34="  # was "92-71="
$t.Cell(5,3).Range.Text = "44+6="  # was "92-10="
$t.Cell(5,4).Range.Text = "3+59="  # was "25+56="
$t.Cell(5,5).Range.Text = "76-72="  # was "33+14="
$t.Cell(6,1).Range.Text = "85-37="  # was "46+48="
$t.Cell(6,2).Range.Text = "33-27="  # was "29-15="
$t.Cell(6,3).Range.Text = "55+3="  # was "43-19="
$t.Cell(6,4).Range.Text = "97-87="  # was "53-53="
$t.Cell(6,5).Range.Text = "13+21="  # was "7+33="
$t.Cell(7,1).Range.Text = "67-61="  # was "13-12="
$t.Cell(7,2).Range.Text = "21-15="  # was "11+17="
$t.Cell(7,3).Range.Text = "98-96="  # was "21+60="
$t.Cell(7,4).Range.Text = "1+87="  # was "67-59="
$t.Cell(7,5).Range.Text = "86-35="  # was "19+51="
$t.Cell(8,1).Range.Text = "43-9="  # was "80-78="
$t.Cell(8,2).Range.Text = "10+32="  # was "78-7="
$t.Cell(8,3).Range.Text = "72-29="  # was "84+9="
$t.Cell(8,4).Range.Text = "55-9="  # was "62-52="
$t.Cell(8,5).Range.Text = "48+14="  # was "73-44="
$t.Cell(9,1).Range.Text = "78-11="  # was "79-31="
$t.Cell(9,2).Range.Text = "46+22="  # was "34-9="
$t.Cell(9,3).Range.Text = "86-57="  # was "44-22="
$t.Cell(9,4).Range.Text = "35-23="  # was "38+48="
$t.Cell(9,5).Range.Text = "13+27="  # was "72-35="
$t.Cell(10,1).Range.Text = "51-48="  # was "53+37="
$t.Cell(10,2).Range.Text = "78+20="  # was "83-15="
$t.Cell(10,3).Range.Text = "78-27="  # was "96-63="
$t.Cell(10,4).Range.Text = "91-46="  # was "92-51="
$t.Cell(10,5).Range.Text = "87-52="  # was "37+0="
$t.Cell(11,1).Range.Text = "79+19="  # was "44+37="
$t.Cell(11,2).Range.Text = "23-0="  # was "43-30="
$t.Cell(11,3).Range.Text = "56-23="  # was "67+4="
$t.Cell(11,4).Range.Text = "57+9="  # was "90-77="
$t.Cell(11,5).Range.Text = "95-83="  # was "32+52="
$t.Cell(12,1).Range.Text = "20-14="  # was "57-37="
$t.Cell(12,2).Range.Text = "53-7="  # was "43+12="
$t.Cell(12,3).Range.Text = "86-43="  # was "50-17="
$t.Cell(12,4).Range.Text = "54-49="  # was "90-69="
$t.Cell(12,5).Range.Text = "55-54="  # was "55-49="
$t.Cell(13,1).Range.Text = "5+8="  # was "63-46="
$t.Cell(13,2).Range.Text = "43-32="  # was "24+58="
$t.Cell(13,3).Range.Text = "74+6="  # was "41-10="
$t.Cell(13,4).Range.Text = "65-56="  # was "55-42="
$t.Cell(13,5).Range.Text = "97-24="  # was "26+41="
$t.Cell(14,1).Range.Text = "17+11="  # was "6+58="
$t.Cell(14,2).Range.Text = "54-29="  # was "16-10="
$t.Cell(14,3).Range.Text = "52+25="  # was "44-33="
$t.Cell(14,4).Range.Text = "9+33="  # was "80+0="
$t.Cell(14,5).Range.Text = "77-21="  # was "24+59="
$t.Cell(15,1).Range.Text = "14+50="  # was "12+43="
$t.Cell(15,2).Range.Text = "58-33="  # was "12+66="
$t.Cell(15,3).Range.Text = "48-11="  # was "2+4="
$t.Cell(15,4).Range.Text = "19+16="  # was "36-31="
$t.Cell(15,5).Range.Text = "92-56="  # was "89-3="
$t.Cell(16,1).Range.Text = "79-73="  # was "4+82="
$t.Cell(16,2).Range.Text = "64-6="  # was "96-38="
$t.Cell(16,3).Range.Text = "8+90="  # was "91-87="
$t.Cell(16,4).Range.Text = "50+36="  # was "75-39="
$t.Cell(16,5).Range.Text = "17+9="  # was "78+5="
$t.Cell(17,1).Range.Text = "38+12="  # was "8-7="
$t.Cell(17,2).Range.Text = "64-27="  # was "74+14="
$t.Cell(17,3).Range.Text = "33+7="  # was "15+42="
$t.Cell(17,4).Range.Text = "55-2="  # was "31+35="
$t.Cell(17,5).Range.Text = "24+26="  # was "90-41="
$t.Cell(18,1).Range.Text = "30+4="  # was "71-65="
$t.Cell(18,2).Range.Text = "49-8="  # was "33+29="
$t.Cell(18,3).Range.Text = "30-22="  # was "21+42="
$t.Cell(18,4).Range.Text = "81+2="  # was "45+54="
$t.Cell(18,5).Range.Text = "57-30="  # was "36+24="
$t.Cell(19,1).Range.Text = "88-9="  # was "22+65="
$t.Cell(19,2).Range.Text = "9+45="  # was "92-62="
$t.Cell(19,3).Range.Text = "9+49="  # was "54-52="
$t.Cell(19,4).Range.Text = "35+29="  # was "79-17="
$t.Cell(19,5).Range.Text = "77-14="  # was "32+20="
$t.Cell(20,1).Range.Text = "50+41="  # was "97-61="
$t.Cell(20,2).Range.Text = "12+54="  # was "47+33="
$t.Cell(20,3).Range.Text = "46+43="  # was "7+92="
$t.Cell(20,4).Range.Text = "32-4="  # was "34+55="
$t.Cell(20,5).Range.Text = "0+2="  # was "12+64="
